$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "VideoYouTube" column (column F, header in shared string 15) was
# removed entirely from the template. Deleting the whole column shifts
# every column after it one position to the left (G->F, H->G, ... P->O),
# and Excel automatically drops the now-unused "VideoYouTube" shared
# string as well as the whole-column data validation that targeted F.
$ws.Columns.Item(6).Delete()

# After the shift, column F ("Informacion adicional", previously column G)
# and column G ("Imagen Logo", previously column H) were resized.
$ws.Columns.Item(6).ColumnWidth = 21.333333333333336
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666

# The selected cell moved to M2 (single cell, no longer the whole row).
[void]$ws.Range("M2").Select()
